# Team-3 Individual Assessment — add "Code Review 2" scoring column (C),
# fix the Code Review 2 total to round like Code Review 1, and update the
# selection/cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Code Review 2 scores (column C) for every student row (2-7),
# mirroring the Code Review 1 (column B) scores.
$ws.Range("C2").Value = 16.6666667
$ws.Range("C3").Value = 16.6666667
$ws.Range("C4").Value = 16.6666667
$ws.Range("C5").Value = 16.6666667
$ws.Range("C6").Value = 16.6666667
$ws.Range("C7").Value = 16.6666667

# The Code Review 2 total should round the same way the Code Review 1
# total already does.
$ws.Range("C8").Formula = "=ROUND(SUM(C2:C7),2)"

# Normalize the formatting on A4 so it matches the rest of the
# Matriculation Number column (centered horizontally + vertically).
$ws.Range("A4").HorizontalAlignment = -4108
$ws.Range("A4").VerticalAlignment = -4108

# Move the active cell/selection to B8 (Code Review 1 total) like it was
# left after the last edit.
$ws.Range("B8").Select() | Out-Null
